$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 31
$ws.Cells(31, 8).Value = 5333.3335
$ws.Cells(31, 9).Value = 4000
$ws.Cells(31, 10).Value = 8000
$ws.Cells(31, 11).Value = 12000
$ws.Cells(31, 12).Value = 24000
$ws.Cells(31, 13).Value = -11770
$ws.Cells(31, 14).Value = -24460

# Row 40
$ws.Cells(40, 8).Value = 1018.375
$ws.Cells(40, 9).Value = 938
$ws.Cells(40, 10).Value = 1054.909
$ws.Cells(40, 11).Value = 938
$ws.Cells(40, 12).Value = 1054.909
$ws.Cells(40, 13).Value = -763
$ws.Cells(40, 14).Value = -1404.909

# Row 52
$ws.Cells(52, 8).Value = 1116.6666
$ws.Cells(52, 10).Value = 1800
$ws.Cells(52, 12).Value = 5400
$ws.Cells(52, 14).Value = -5720

# Row 127
$ws.Cells(127, 8).Value = 1047.75
$ws.Cells(127, 9).Value = 691.125
$ws.Cells(127, 10).Value = 1226.0625
$ws.Cells(127, 11).Value = 2073.375
$ws.Cells(127, 12).Value = 3678.1875
$ws.Cells(127, 13).Value = 2886.625
$ws.Cells(127, 14).Value = -13598.1875

$ws = $wb.Worksheets.Item("ARM")
# Row 26
$ws.Cells(26, 8).Value = 1782.0714
$ws.Cells(26, 9).Value = 1534.5385
$ws.Cells(26, 10).Value = 5000
$ws.Cells(26, 11).Value = 1534.5385
$ws.Cells(26, 12).Value = 5000
$ws.Cells(26, 13).Value = -1204.5385
$ws.Cells(26, 14).Value = -5660

# Row 27
$ws.Cells(27, 8).Value = 2000
$ws.Cells(27, 9).Value = 2000
$ws.Cells(27, 11).Value = 2000
$ws.Cells(27, 13).Value = -1816

# Row 32
$ws.Cells(32, 8).Value = 2431.88
$ws.Cells(32, 9).Value = 2022.7906
$ws.Cells(32, 11).Value = 2022.7906
$ws.Cells(32, 13).Value = -1735.7906

# Row 34
$ws.Cells(34, 8).Value = 0
$ws.Cells(34, 10).Value = 0
$ws.Cells(34, 12).Value = 0
$ws.Cells(34, 14).ClearContents()

# Row 40
$ws.Cells(40, 8).Value = 0
$ws.Cells(40, 10).Value = 0
$ws.Cells(40, 12).Value = 0
$ws.Cells(40, 14).ClearContents()

# Row 45
$ws.Cells(45, 8).Value = 1140.7368
$ws.Cells(45, 9).Value = 928.2222
$ws.Cells(45, 10).Value = 1332
$ws.Cells(45, 11).Value = 928.2222
$ws.Cells(45, 12).Value = 1332
$ws.Cells(45, 13).Value = -551.2222
$ws.Cells(45, 14).Value = -2086

$ws = $wb.Worksheets.Item("BSM")
# Row 16
$ws.Cells(16, 8).Value = 0
$ws.Cells(16, 10).Value = 0
$ws.Cells(16, 12).Value = 0
$ws.Cells(16, 14).ClearContents()

# Row 34
$ws.Cells(34, 8).Value = 10000
$ws.Cells(34, 10).Value = 10000
$ws.Cells(34, 12).Value = 10000
$ws.Cells(34, 14).Value = -10228

# Row 86
$ws.Cells(86, 8).Value = 9120.429
$ws.Cells(86, 9).Value = 10384.333
$ws.Cells(86, 10).Value = 8172.5
$ws.Cells(86, 11).Value = 10384.333
$ws.Cells(86, 12).Value = 8172.5
$ws.Cells(86, 13).Value = -9261.333000000001
$ws.Cells(86, 14).Value = -10418.5

# Row 89
$ws.Cells(89, 8).Value = 9120.429
$ws.Cells(89, 9).Value = 10384.333
$ws.Cells(89, 10).Value = 8172.5
$ws.Cells(89, 11).Value = 51921.665
$ws.Cells(89, 12).Value = 40862.5
$ws.Cells(89, 13).Value = -46305.665
$ws.Cells(89, 14).Value = -52094.5

$ws = $wb.Worksheets.Item("CRP")
# Row 14
$ws.Cells(14, 8).Value = 997.5
$ws.Cells(14, 9).Value = 830
$ws.Cells(14, 11).Value = 830
$ws.Cells(14, 13).Value = -660

# Row 17
$ws.Cells(17, 9).Value = 3000
$ws.Cells(17, 11).Value = 3000
$ws.Cells(17, 13).Value = -2826

# Row 25
$ws.Cells(25, 8).Value = 10000
$ws.Cells(25, 9).Value = 0
$ws.Cells(25, 11).Value = 0
$ws.Cells(25, 13).ClearContents()

# Row 57
$ws.Cells(57, 8).Value = 21000
$ws.Cells(57, 10).Value = 21000
$ws.Cells(57, 12).Value = 21000
$ws.Cells(57, 14).Value = -22120

# Row 135
$ws.Cells(135, 8).Value = 36195
$ws.Cells(135, 10).Value = 27097.5
$ws.Cells(135, 12).Value = 27097.5
$ws.Cells(135, 14).Value = -37237.5

$ws = $wb.Worksheets.Item("CUL")
# Row 39
$ws.Cells(39, 8).Value = 2410.6428
$ws.Cells(39, 9).Value = 0
$ws.Cells(39, 10).Value = 2410.6428
$ws.Cells(39, 11).Value = 0
$ws.Cells(39, 12).Value = 7231.928400000001
$ws.Cells(39, 13).ClearContents()
$ws.Cells(39, 14).Value = -7819.928400000001

# Row 40
$ws.Cells(40, 8).Value = 187.30435
$ws.Cells(40, 9).Value = 105.68421
$ws.Cells(40, 10).Value = 575
$ws.Cells(40, 11).Value = 422.73684
$ws.Cells(40, 12).Value = 2300
$ws.Cells(40, 13).Value = -353.73684
$ws.Cells(40, 14).Value = -2438

# Row 41
$ws.Cells(41, 8).Value = 616
$ws.Cells(41, 9).Value = 360
$ws.Cells(41, 10).Value = 1000
$ws.Cells(41, 11).Value = 1080
$ws.Cells(41, 12).Value = 3000
$ws.Cells(41, 13).Value = -742
$ws.Cells(41, 14).Value = -3676

# Row 46
$ws.Cells(46, 8).Value = 1344.4445
$ws.Cells(46, 9).Value = 650
$ws.Cells(46, 10).Value = 1542.8572
$ws.Cells(46, 11).Value = 1950
$ws.Cells(46, 12).Value = 4628.571599999999
$ws.Cells(46, 13).Value = -1859
$ws.Cells(46, 14).Value = -4810.571599999999

$ws = $wb.Worksheets.Item("GSM")
# Row 15
$ws.Cells(15, 8).Value = 8400
$ws.Cells(15, 10).Value = 8400
$ws.Cells(15, 12).Value = 8400
$ws.Cells(15, 14).Value = -8976

# Row 23
$ws.Cells(23, 8).Value = 8333.333000000001
$ws.Cells(23, 9).Value = 0
$ws.Cells(23, 10).Value = 8333.333000000001
$ws.Cells(23, 11).Value = 0
$ws.Cells(23, 12).Value = 8333.333000000001
$ws.Cells(23, 13).ClearContents()
$ws.Cells(23, 14).Value = -8779.333000000001

# Row 40
$ws.Cells(40, 8).Value = 0
$ws.Cells(40, 10).Value = 0
$ws.Cells(40, 12).Value = 0
$ws.Cells(40, 14).ClearContents()

# Row 46
$ws.Cells(46, 8).Value = 4500
$ws.Cells(46, 9).Value = 4500
$ws.Cells(46, 11).Value = 4500
$ws.Cells(46, 13).Value = -4344

# Row 81
$ws.Cells(81, 8).Value = 8400
$ws.Cells(81, 10).Value = 8400
$ws.Cells(81, 12).Value = 8400
$ws.Cells(81, 14).Value = -10396

# Row 84
$ws.Cells(84, 8).Value = 8400
$ws.Cells(84, 10).Value = 8400
$ws.Cells(84, 12).Value = 25200
$ws.Cells(84, 14).Value = -35184

# Row 102
$ws.Cells(102, 8).Value = 3427.7878
$ws.Cells(102, 9).Value = 3607.1785
$ws.Cells(102, 11).Value = 3607.1785
$ws.Cells(102, 13).Value = -1985.1785

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Cells(22, 8).Value = 748
$ws.Cells(22, 10).Value = 748
$ws.Cells(22, 12).Value = 748
$ws.Cells(22, 14).Value = -1338

# Row 27
$ws.Cells(27, 8).Value = 748
$ws.Cells(27, 10).Value = 748
$ws.Cells(27, 12).Value = 748
$ws.Cells(27, 14).Value = -962

# Row 47
$ws.Cells(47, 8).Value = 24995.455
$ws.Cells(47, 10).Value = 24995.455
$ws.Cells(47, 12).Value = 24995.455
$ws.Cells(47, 14).Value = -25975.455

# Row 52
$ws.Cells(52, 8).Value = 24995.455
$ws.Cells(52, 10).Value = 24995.455
$ws.Cells(52, 12).Value = 24995.455
$ws.Cells(52, 14).Value = -25461.455

# Row 58
$ws.Cells(58, 8).Value = 2250
$ws.Cells(58, 9).Value = 2250
$ws.Cells(58, 11).Value = 2250
$ws.Cells(58, 13).Value = -1990

# Row 61
$ws.Cells(61, 8).Value = 1714.1538
$ws.Cells(61, 9).Value = 1498.5454
$ws.Cells(61, 10).Value = 2900
$ws.Cells(61, 11).Value = 1498.5454
$ws.Cells(61, 12).Value = 2900
$ws.Cells(61, 13).Value = -1296.5454
$ws.Cells(61, 14).Value = -3304

# Row 113
$ws.Cells(113, 8).Value = 1714.1538
$ws.Cells(113, 9).Value = 1498.5454
$ws.Cells(113, 10).Value = 2900
$ws.Cells(113, 11).Value = 1498.5454
$ws.Cells(113, 12).Value = 2900
$ws.Cells(113, 13).Value = 671.4546
$ws.Cells(113, 14).Value = -7240

# Row 136
$ws.Cells(136, 8).Value = 5208.108
$ws.Cells(136, 9).Value = 2940.8333
$ws.Cells(136, 10).Value = 7356.0527
$ws.Cells(136, 11).Value = 8822.499899999999
$ws.Cells(136, 12).Value = 22068.1581
$ws.Cells(136, 13).Value = -6272.499899999999
$ws.Cells(136, 14).Value = -27168.1581

$ws = $wb.Worksheets.Item("WVR")
# Row 28
$ws.Cells(28, 8).Value = 4151.5
$ws.Cells(28, 9).Value = 2995
$ws.Cells(28, 10).Value = 4729.75
$ws.Cells(28, 11).Value = 2995
$ws.Cells(28, 12).Value = 4729.75
$ws.Cells(28, 13).Value = -2647
$ws.Cells(28, 14).Value = -5425.75

# Row 40
$ws.Cells(40, 8).Value = 9000
$ws.Cells(40, 10).Value = 9000
$ws.Cells(40, 12).Value = 9000
$ws.Cells(40, 14).Value = -9298
